# Update "Förändrad" (column C) date value from 45175 (2023-09-06) to
# 45177 (2023-09-08) for every data row (rows 2 through 480).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C480").Value = 45177
